# Daily attendance processing - 2026-01-25 10:35:34
# Normalizes the "Recorded By" column (G) so that entries of the form
# "dnasr281@gmail.com, <other>" have their two parts swapped to
# "<other>, dnasr281@gmail.com".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

$prefix = "dnasr281@gmail.com, "

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # Column G
    $val = $cell.Value2

    if ($val -ne $null -and $val -like "$prefix*") {
        $other = $val.Substring($prefix.Length)
        $cell.Value2 = "$other, dnasr281@gmail.com"
    }
}
